$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160:176 down to 161:177.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new weekly price record.
$ws.Range("A160").Value = 4
$ws.Range("B160").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C160").Value = "Los Lagos"
$ws.Range("D160").Value = 45212
$ws.Range("D160").NumberFormat = $ws.Range("D161").NumberFormat
$ws.Range("E160").Value = 10
$ws.Range("F160").Value = 100112031
$ws.Range("G160").Value = "Poroto verde"
$ws.Range("H160").Value = "Magnum"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 40
$ws.Range("K160").Value = 37000
$ws.Range("L160").Value = 37000
$ws.Range("M160").Value = 37000
$ws.Range("N160").Value = "$/malla 25 kilos"
$ws.Range("O160").Value = "Perú"
$ws.Range("P160").Value = 1480
$ws.Range("Q160").Value = 25
$ws.Range("R160").Value = "Hortaliza"
